# Update values on Delivery_results sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Delivery_results")

$ws1.Range("D2").Value = 22
$ws1.Range("E2").Value = 35.48
$ws1.Range("F2").Value = 2.92
$ws1.Range("G2").Value = 5.59

$ws1.Range("D3").Value = 15
$ws1.Range("E3").Value = 25
$ws1.Range("F3").Value = 4.5
$ws1.Range("G3").Value = 6.8

$ws1.Range("F4").Value = 1.93
$ws1.Range("G4").Value = 4.76

$ws1.Range("F5").Value = 3.43
$ws1.Range("G5").Value = 5.38

$ws1.Range("D6").Value = 24
$ws1.Range("E6").Value = 64.86
$ws1.Range("F6").Value = 2.92
$ws1.Range("G6").Value = 6.75

# Add new sheet "Total_distance" after Total_emissions
$ws2 = $wb.Worksheets.Item("Total_emissions")
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Total_distance"

$ws3.Range("A1").Value = "Total distance (km)"
$ws2.Range("A1").Copy()
$ws3.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws3.Range("A2").Value = 39.76204
